# The "Test Data" cell for the "missing parentheses" addCandidates test
# contained two paragraphs:
#   candidatesLine: "C0 (P0), C1 P1"
#   line: 3
# The stray "line: 3" paragraph (which doesn't carry useful parsing data)
# is removed so the cell is left with just the candidatesLine paragraph.

$d = $word.ActiveDocument

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "line: 3") {
        $p.Range.Delete()
    }
}
